$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$rng = $ws.Range("A271")
try {
  $s = $wb.Styles.Add("MyStyle")
  Write-Host "added style"
} catch {
  Write-Host "ERR: $_"
}
